$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted before the current row 44;
# every following record (rows 44-95) shifts down one row to 45-96.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new record.
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 44942
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100112030
$ws.Range("G44").Value = "Poroto granado"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 45
$ws.Range("K44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("M44").Value = 50000
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Región del Maule"
$ws.Range("P44").Value = 2000
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
